$d = $word.ActiveDocument

# Locate the paragraph that holds the site footer / copyright notice
# ("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll ...").
$target = $null
foreach ($p in @($d.Paragraphs)) {
    if ($p.Range.Text -like "*Contact: luizeleno@usp.br*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # The paragraph immediately before the footer is an empty paragraph
    # that carries the page-break-before formatting for the footer page;
    # it goes away together with the footer text. The paragraph
    # immediately after the footer is also an empty "spacer" paragraph
    # that is removed along with it, so the range to delete runs from
    # the start of the preceding paragraph through the end of the
    # paragraph mark of the following (spacer) paragraph.
    $prev = $target.Previous(1)
    $next = $target.Next(1)

    $delStart = $prev.Range.Start
    $delEnd = $next.Range.End

    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
